$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# --- Sheet 1: Creacion de canal de Youtube ---
# Punteo (grades) for activity rows 8-12 (Punteo column E).
# These cells are locked on a protected sheet, so we briefly unlock them to
# write the value, then restore their exact original formatting (style s="1")
# by pasting formats from an untouched sibling cell (E7) - this avoids
# calling Unprotect/Protect, which would otherwise strip the sheetProtection
# flags (objects/scenarios/selectLockedCells) from the saved file.
$ws1 = $wb.Worksheets.Item("Creacion de canal de Youtube")

foreach ($row in 8..12) {
    $ws1.Range("E$row").Locked = $false
    $ws1.Range("E$row").Value = 5
}
$ws1.Range("E7").Copy()
$ws1.Range("E8:E12").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# --- Sheet 2: Explicacion de Algoritmos ---
# Responsable (names) and Punteo (grades) for activity rows 7-11.
$ws2 = $wb.Worksheets.Item("Explicacion de Algoritmos")

# Column D cells are already unlocked on the protected sheet, so values can
# be written directly.
$ws2.Range("D7").Value = "Erick "
$ws2.Range("D8").Value = "Kevin"
$ws2.Range("D9").Value = "Juan"
$ws2.Range("D10").Value = "Luis"
$ws2.Range("D11").Value = "Julio"

# Column E cells are locked; use the same unlock/write/restore-format trick.
foreach ($row in 7..11) {
    $ws2.Range("E$row").Locked = $false
}
$ws2.Range("E7").Value = 5
$ws2.Range("E8").Value = 4
$ws2.Range("E9").Value = 5
$ws2.Range("E10").Value = 4
$ws2.Range("E11").Value = 4

$ws2.Range("E12").Copy()
$ws2.Range("E7:E11").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

# --- Update saved selections to match the final view state ---
# Select sheet2's cell first, then sheet1's cell last so sheet1 (the
# originally active tab) ends up active again in the saved file.
$ws2.Range("D13").Select()
$ws1.Range("C14").Select()
